$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.777.82"
$ws.Range("E2").Value = "  -3.07%  "

$ws.Range("D3").Value = "3.436.50"
$ws.Range("E3").Value = "  -2.78%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'570.40"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").Value = "'174.54"
$ws.Range("E6").Value = "  -7.43%  "

$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  -1.55%  "

$ws.Range("E10").Value = "  +4.87%  "

$ws.Range("D11").Value = "'55.00"
$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").Value = "'0.0000273"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").Value = "'9.11"
$ws.Range("E13").Value = "  -3.10%  "

$ws.Range("D14").Value = "3.984.95"
$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.121"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.430.07"
$ws.Range("E16").Value = "  -2.99%  "

$ws.Range("D17").Value = "'18.10"
$ws.Range("E17").Value = "  -0.85%  "

$ws.Range("D18").Value = "'11.86"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").Value = "64.774.62"
$ws.Range("E19").Value = "  -3.18%  "

$ws.Range("D20").Value = "'0.989"
$ws.Range("E20").Value = "  -0.82%  "

$ws.Range("D21").Value = "'406.86"
$ws.Range("E21").Value = "  -4.97%  "

$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("E23").Value = "  +7.90%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'13.38"
$ws.Range("E24").Value = "  +9.14%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'83.60"
$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("D26").Value = "'10.80"
$ws.Range("E26").Value = "  -3.14%  "

$ws.Range("D27").Value = "'2.80"
$ws.Range("E27").Value = "  -3.29%  "

$ws.Range("D28").Value = "'9.01"
$ws.Range("E28").Value = "  -2.56%  "

$ws.Range("E29").Value = "  -1.56%  "

$ws.Range("D30").Value = "'6.55"
$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("D31").Value = "'11.53"
$ws.Range("E31").Value = "  -1.81%  "

$ws.Range("D32").Value = "'586.35"
$ws.Range("E32").Value = "  -8.77%  "

$ws.Range("E33").Value = "  -3.25%  "

$ws.Range("D34").Value = "'59.71"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("E35").Value = "  +2.85%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'3.56"
$ws.Range("E37").Value = "  +6.11%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0775"
$ws.Range("E38").Value = "  -4.76%  "

$ws.Range("D39").Value = "'36.17"
$ws.Range("E39").Value = "  -6.31%  "

$ws.Range("D40").Value = "'0.376"
$ws.Range("E40").Value = "  -4.26%  "

$ws.Range("D41").Value = "3.183.70"
$ws.Range("E41").Value = "  +4.35%  "

$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("E43").Value = "  +1.20%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'3.24"
$ws.Range("E44").Value = "  -3.97%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.50"
$ws.Range("E45").Value = "  -6.44%  "

$ws.Range("D46").Value = "'0.0409"
$ws.Range("E46").Value = "  -2.78%  "

$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("E48").Value = "  -5.39%  "

$ws.Range("D49").Value = "'8.43"
$ws.Range("E49").Value = "  -2.55%  "

$ws.Range("D50").Value = "'136.64"
$ws.Range("E50").Value = "  -4.74%  "

$ws.Range("D51").Value = "'2.32"
